$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two employee names in column C (rows 2 and 3).
$ws.Range("C2").Value = "Ranga  Akunuri"
$ws.Range("C3").Value = "Peter Mac Anderson"

# Give C2:C3 their own distinct style: same font as the header cells
# (Arial 9.6 FF374151) but vertically centered / wrapped like the other
# body cells. Copying the header's format keeps the existing font entry
# (no duplicate font created); the vertical alignment tweak after the
# paste creates a new cellXf distinct from the plain body style.
$ws.Range("B1").Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)
$ws.Range("C2:C3").VerticalAlignment = -4108

# Row 3 grows a bit taller to fit the longer wrapped name.
$ws.Rows.Item(3).RowHeight = 38.25

# Move the active selection.
$ws.Range("F8").Select()
